$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = $false
$ws.Range("B1").Value = "simple"
$ws.Range("C1").Value = 7.112169692521646

$ws.Range("A2").Value = $true
$ws.Range("B2").Value = "simple"
$ws.Range("C2").Value = 5.315955288750741

$ws.Range("A3").Value = $false
$ws.Range("B3").Value = "hyper_heuristic"
$ws.Range("C3").Value = 5.315955288750741

$ws.Range("A4").Value = $true
$ws.Range("B4").Value = "hyper_heuristic"
$ws.Range("C4").Value = 5.315955288750741

$ws.Range("A1:C2").Select() | Out-Null
